$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 473
$ws.Range("F5").Value = 2488
$ws.Range("F7").Value = 81
$ws.Range("F8").Value = 89
$ws.Range("D9").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("F9").Value = 1713
$ws.Range("D10").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("F10").Value = 1713
$ws.Range("F12").Value = 80
$ws.Range("F14").Value = 24
$ws.Range("F16").Value = 952
$ws.Range("F17").Value = 339
$ws.Range("F20").Value = 7584
$ws.Range("F21").Value = 7584
$ws.Range("F22").Value = 8630
$ws.Range("F23").Value = 60
$ws.Range("F27").Value = 102
$ws.Range("F29").Value = 23
$ws.Range("D33").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("F33").Value = 1521
$ws.Range("F34").Value = 23
$ws.Range("F38").Value = 304
$ws.Range("F40").Value = 811
$ws.Range("D42").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("F42").Value = 1377
$ws.Range("F46").Value = 99
$ws.Range("F48").Value = 6
$ws.Range("D49").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("F49").Value = 201
$ws.Range("D50").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 39
$ws.Range("F15").Value = 25
$ws.Range("F19").Value = 4
$ws.Range("F20").Value = 319
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2665
$ws.Range("G3").Value = 0
$ws.Range("F4").Value = 310
$ws.Range("F5").Value = 160
$ws.Range("F6").Value = 23
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 473
$ws.Range("F5").Value = 2665
$ws.Range("G5").Value = 30
$ws.Range("F6").Value = 310
$ws.Range("F7").Value = 160
$ws.Range("F8").Value = 39
$ws.Range("F9").Value = 2488
$ws.Range("F11").Value = 89
$ws.Range("D12").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("F12").Value = 1713
$ws.Range("D13").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("F13").Value = 1713
$ws.Range("F15").Value = 24
$ws.Range("F17").Value = 952
$ws.Range("F18").Value = 339
$ws.Range("F23").Value = 7584
$ws.Range("F24").Value = 8630
$ws.Range("F25").Value = 60
$ws.Range("F27").Value = 102
$ws.Range("B29").Value = "2024-10-02"
$ws.Range("C29").Value = "北京·第19届IJOY漫展【217专场见面会】"
$ws.Range("D29").Value = "天辰东路7号 北京国家会议中心"
$ws.Range("E29").Value = "2024.10.02 12:25-10.02 16:30"
$ws.Range("F29").Value = 23
$ws.Range("G29").Value = 168
$ws.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=91561"
$ws.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202408/cAghXlck1724928163645.jpeg"
$ws.Range("B30").Value = "2024-10-03"
$ws.Range("C30").Value = "北京·明日方舟同人only-厮守序言"
$ws.Range("D30").Value = "安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园"
$ws.Range("E30").Value = "2024.10.03 09:30-10.03 17:00"
$ws.Range("F30").Value = 362
$ws.Range("G30").Value = 60
$ws.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=90959"
$ws.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202408/rIGY6eyZ1723974119991.jpeg"
$ws.Range("C31").Value = "北京·第五人格only同人展"
$ws.Range("D31").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("E31").Value = "2024.10.04 10:00-10.04 17:00"
$ws.Range("F31").Value = 1521
$ws.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=89309"
$ws.Range("I31").Value = "//i0.hdslb.com/bfs/openplatform/202407/4XsICpa71721046044404.jpeg"
$ws.Range("B32").Value = "2024-10-04"
$ws.Range("C32").Value = "帝都·重返未来1999同人ONLY金秋深眠"
$ws.Range("D32").Value = "华佗路与新源大街交汇处西100米 凯德MALL·大兴"
$ws.Range("E32").Value = "2024.10.04 10:00-10.05 17:00"
$ws.Range("F32").Value = 23
$ws.Range("G32").Value = 68
$ws.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=92315"
$ws.Range("I32").Value = "//i1.hdslb.com/bfs/openplatform/202409/YHMYHehz1726129707544.jpeg"
$ws.Range("F35").Value = 304
$ws.Range("F38").Value = 811
$ws.Range("D41").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("F41").Value = 1377
$ws.Range("F45").Value = 99
$ws.Range("F47").Value = 6
$ws.Range("D48").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"
$ws.Range("F48").Value = 201
$ws.Range("F50").Value = 319
$ws.Range("D51").Value = "北花园路1号超级蜂巢C座 超级蜂巢国际会议中心"

Write-Host "Applied 98 cell edits"